# Foi adicionado um produto
# Add a new product ("Jodo da tia vera") to ESTOQUE, with a matching
# ENTRADA (stock-in) record of 10 units from Alciomar, and refresh the
# PESSOA roll-up that depends on it.
#
# NOTE: this workbook has a defined name `ESTOQUE` (ESTOQUE!$A$2:$F$74)
# that overlaps the full-column defined name `COD_ESTOQUE`
# (ESTOQUE!$A:$A) on the same sheet/column. A full automatic recalc after
# writing into ESTOQUE/ENTRADA column A corrupts unrelated VLOOKUP(...,
# ESTOQUE, ...) lookups elsewhere (ENTRADA/SAIDA) in this runtime. Work
# around it by switching to manual calculation for the duration of the
# edit: every *newly entered* formula still gets evaluated on entry (so
# our new cells compute correctly), but already-evaluated formula cells
# are left with their last-good cached value instead of being corrupted
# by a full-sheet recalc. We only re-enter the couple of PESSOA formulas
# that genuinely must change, which forces just those to recompute.

$excel.Calculation = -4135   # xlCalculationManual

$wb = $excel.ActiveWorkbook

# ---- ESTOQUE: add the new product row (table auto-expands) ----------
$wsEstoque = $wb.Worksheets.Item("ESTOQUE")
$tabEstoque = $wsEstoque.ListObjects.Item("Tabela4")
$tabEstoque.ListRows.Add() | Out-Null

$wsEstoque.Range("A4").Value = 3
$wsEstoque.Range("B4").Value = "Jodo da tia vera"
$wsEstoque.Range("C4").Value = 15

# ---- ENTRADA: add the stock-in record for the new product -----------
$wsEntrada = $wb.Worksheets.Item("ENTRADA")
$tabEntrada = $wsEntrada.ListObjects.Item("Tabela2")
$tabEntrada.ListRows.Add() | Out-Null

$wsEntrada.Range("A6").Value = 3
$wsEntrada.Range("B6").Formula = "=IF(A6=`"`",`"`",VLOOKUP(A6,ESTOQUE,2,0))"
$wsEntrada.Range("C6").Formula = "=IF(A6=`"`",`"`",VLOOKUP(A6,ESTOQUE,3,0))"
$wsEntrada.Range("D6").Value = "5/7/2017"
$wsEntrada.Range("E6").Value = 10
$wsEntrada.Range("F6").Value = "Alciomar"

# ---- ESTOQUE: fill in the calculated columns now that ENTRADA has the
# new row in place, so they pick up the right totals on first entry ----
$wsEstoque.Range("D4").Formula = "=IF(A4=`"`",`"`",SUMIF(COD_ENTRADA,COD_ESTOQUE,COD_QUANTIDADE))"
$wsEstoque.Range("E4").Formula = "=IF(A4=`"`",`"`",SUMIF(COD_SAIDA,COD_ESTOQUE,COD_QUANTIDADE_SAIDA))"
$wsEstoque.Range("F4").Formula = "=IF(D4=`"`",`"`",D4-E4)"

# ---- PESSOA: Alciomar's ENTRADA/TOTAL roll-up must reflect the new ---
# ---- stock-in record; re-enter the formulas to force a refresh ------
$wsPessoa = $wb.Worksheets.Item("PESSOA")
$wsPessoa.Range("E2").Formula = "=SUMIF(ENTRADA!F:F,T_PESSOA[[#This Row],[NOME]],COD_QUANTIDADE)"
$wsPessoa.Range("G2").Formula = "=T_PESSOA[[#This Row],[ENTRADA]]-T_PESSOA[[#This Row],[SAIDA]]"

# ---- Selections / active sheet, matching the edit session -----------
$wsEstoque.Activate()
$wsEstoque.Range("C5").Select() | Out-Null

$wsEntrada.Activate()
$wsEntrada.Range("F6").Select() | Out-Null

$wsPessoa.Activate()

$excel.Calculation = -4105   # xlCalculationAutomatic
